# "Generate Report for Handback"
# For each locale sheet (zh-cn, de-de) and each tracked source file (rows 2 & 3):
#   - Status moves from "Ready for handoff" -> "Handed back: in sync with en-US"
#     (this text is shared by the Overview sheet's per-locale status columns too)
#   - Latest Target File (col I) is filled in with the source file name, hyperlinked
#     to the same URL as the Source File Name link in col A
#   - Latest Handback File (col J) is filled in with the generated xlf file name
#   - Latest Handback DateTime (col K) is filled in with the handback timestamp
#     (zh-cn processed first, de-de processed ~17s later)

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$file1 = "b64778cb-8792-415b-80c6-c326caff7005.md"
$file2 = "ff3b416d-95d1-4df8-aca0-c7d146132f98.md"
$file1Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e1708cea0d9e0c69516bfc17e7db7d3915b2623/e2e/$file1"
$file2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e1708cea0d9e0c69516bfc17e7db7d3915b2623/e2e/$file2"

# ---------------------------------------------------------------------------
# Overview sheet: both locale-status columns for both files pick up the new
# status text (it's the same shared label used on the per-locale sheets).
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# Per-locale sheets: zh-cn and de-de
# ---------------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-09-07 03:26:22";
       Xlf1 = "b64778cb-8792-415b-80c6-c326caff7005.692fb0777a720d9af51e0cc0a842fb4369c610c1.zh-cn.xlf";
       Xlf2 = "ff3b416d-95d1-4df8-aca0-c7d146132f98.b1d92737a1b7220bcdc23947b9e5ffdd2d6804a6.zh-cn.xlf" },
    @{ Name = "de-de"; HandbackTime = "2016-09-07 03:26:39";
       Xlf1 = "b64778cb-8792-415b-80c6-c326caff7005.692fb0777a720d9af51e0cc0a842fb4369c610c1.de-de.xlf";
       Xlf2 = "ff3b416d-95d1-4df8-aca0-c7d146132f98.b1d92737a1b7220bcdc23947b9e5ffdd2d6804a6.de-de.xlf" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    # Status column
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File -> hyperlink to the source .md file (same target as col A)
    $ws.Hyperlinks.Add($ws.Range("I2"), $file1Url, $null, $null, $file1)
    $ws.Hyperlinks.Add($ws.Range("I3"), $file2Url, $null, $null, $file2)

    # Latest Handback File
    $ws.Range("J2").Value = $locale.Xlf1
    $ws.Range("J3").Value = $locale.Xlf2

    # Latest Handback DateTime
    $ws.Range("K2").Value = $locale.HandbackTime
    $ws.Range("K3").Value = $locale.HandbackTime

    # Column widths widened to fit the newly-populated columns
    $ws.Columns.Item(3).ColumnWidth = 29.1
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}
